$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosCuenta")

$ws.Range("A2").Value = "PMSmoke"
$ws.Range("B2").Value = "ApellidoPMSmoke"
$ws.Range("C2").Value = 27100105
$ws.Range("D2").Value = 108
